# Rotate columns E (group-code), F (group-name), G (category-name) to the
# right for every row, including the header row:
#   new E = old G   (category-name moves into the group-code slot)
#   new F = old E   (group-code moves into the group-name slot)
#   new G = old F   (group-name moves into the category-name slot)
#
# Net effect: the header becomes
#   code | name | status | category-code | category-name | group-code | group-name
# i.e. category-code/category-name are adjacent, and group-code/group-name are
# adjacent (matching the published codeforiati codelist re-export).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

$eVals = New-Object 'object[]' $lastRow
$fVals = New-Object 'object[]' $lastRow
$gVals = New-Object 'object[]' $lastRow

for ($r = 1; $r -le $lastRow; $r++) {
    $eVals[$r - 1] = $ws.Cells.Item($r, 5).Value()
    $fVals[$r - 1] = $ws.Cells.Item($r, 6).Value()
    $gVals[$r - 1] = $ws.Cells.Item($r, 7).Value()
}

# Column F receives the old group-code values, which are digit-only strings
# (e.g. "110"). Force the column to text first so Excel doesn't silently
# reinterpret them as numbers, then restore the default style so the cells
# end up identical (no explicit style) to the rest of the sheet.
$ws.Columns.Item(6).NumberFormat = "@"

for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = $gVals[$r - 1]
    $ws.Cells.Item($r, 6).Value = $eVals[$r - 1]
    $ws.Cells.Item($r, 7).Value = $fVals[$r - 1]
}

$ws.Columns.Item(6).Style = "Normal"
